$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 112394378
$ws.Range("B3").Value = 89571
$ws.Range("Q3").Value = 749921
$ws.Range("R3").Value = 7535985

# Row 4
$ws.Range("A4").Value = 112394252
$ws.Range("B4").Value = 89571
$ws.Range("Q4").Value = 749928
$ws.Range("R4").Value = 7535991

# Row 5
$ws.Range("B5").Value = 89834

# Row 6
$ws.Range("A6").Value = 112438309
$ws.Range("B6").Value = 90448
$ws.Range("E6").Value = 4745
$ws.Range("F6").Value = "Tallriska"
$ws.Range("G6").Value = "Lactarius musteus"
$ws.Range("H6").Value = "Fr."
$ws.Range("Q6").Value = 749892
$ws.Range("R6").Value = 7535980
$ws.Range("S6").Value = 25
$ws.Range("Z6").Value = "11:54"
$ws.Range("AB6").Value = "11:54"
$ws.Range("AX6").Value = "Stefan Andersson, Christina Boyd, per-erik mukka"

# Row 7
$ws.Range("A7").Value = 112457286
$ws.Range("B7").Value = 89834
$ws.Range("E7").Value = 658
$ws.Range("F7").Value = "Rosenticka"
$ws.Range("G7").Value = "Rhodofomes roseus"
$ws.Range("H7").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q7").Value = 749868
$ws.Range("R7").Value = 7536000
$ws.Range("S7").Value = 10
$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()
$ws.Range("AX7").Value = "Stefan Andersson, per-erik mukka, Christina Boyd"

# Row 8
$ws.Range("A8").Value = 112501187
$ws.Range("B8").Value = 89903
$ws.Range("D8").Value = "VU"
$ws.Range("E8").Value = 1506
$ws.Range("F8").Value = "Ostticka"
$ws.Range("G8").Value = "Skeletocutis odora"
$ws.Range("H8").Value = "(Sacc.) Ginns"

# Row 9
$ws.Range("A9").Value = 112501192
$ws.Range("Q9").Value = 749922
$ws.Range("R9").Value = 7535992

# Row 10
$ws.Range("A10").Value = 112501206
$ws.Range("B10").Value = 89553
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 1202
$ws.Range("F10").Value = "Ullticka"
$ws.Range("G10").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H10").Value = "(P.Karst.) Fiasson & Niemelä"

# Row 11
$ws.Range("A11").Value = 112501052
$ws.Range("Q11").Value = 749927
$ws.Range("R11").Value = 7535984

# Row 12
$ws.Range("A12").Value = 112501198
$ws.Range("B12").Value = 89834
$ws.Range("E12").Value = 658
$ws.Range("F12").Value = "Rosenticka"
$ws.Range("G12").Value = "Rhodofomes roseus"
$ws.Range("H12").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
